$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new trade row (row 10), matching formatting of row 9 above it
$ws.Range("A9:I9").Copy()
$ws.Range("A10:I10").PasteSpecial(-4122)
$ws.Range("A10:I10").PasteSpecial(-4163)

$ws.Cells.Item(10, 1).Value = 42654.745983796296
$ws.Cells.Item(10, 2).Value = $true
$ws.Cells.Item(10, 3).Value = 10022.16
$ws.Cells.Item(10, 4).Value = 10009.15
$ws.Cells.Item(10, 5).Value = 18.870000999999998
$ws.Cells.Item(10, 6).Value = 18.920000000000002
$ws.Cells.Item(10, 7).Value = $false
$ws.Cells.Item(10, 8).Value = 0.26
$ws.Cells.Item(10, 9).Value = $false

# Auto-fit columns to best fit like the repeater tool does
$ws.Columns("A:I").AutoFit()

# Ensure the resulting column widths match Excel's own "best fit" pixel metrics
$ws.Columns("A").ColumnWidth = 14.5
$ws.Columns("B").ColumnWidth = 7.333333333333333
$ws.Columns("C").ColumnWidth = 8.0
$ws.Columns("D").ColumnWidth = 10.333333333333334
$ws.Columns("E").ColumnWidth = 9.0
$ws.Columns("F").ColumnWidth = 6.166666666666667
$ws.Columns("G").ColumnWidth = 9.5
$ws.Columns("H").ColumnWidth = 13.833333333333334
$ws.Columns("I").ColumnWidth = 11.0
